$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (MBDC1731 / COLOMBOS) ---
$ws.Range("E9").Value = $false
$ws.Range("F9").Value = "DOES NOT LOAD"
# G9 gets a date-style number format applied before the text value is entered,
# matching Excel's "d-mmm-yy" built-in format (numFmtId 15) while the stored
# content itself is text.
$ws.Range("G9").NumberFormat = "d-mmm-yy"
$ws.Range("G9").Value = "march 21 2022"
$ws.Range("H9").Value = "tried using multiple different browsers and devices, website would not load"

# --- Row 10 (MBDC1405 / COMBREX) ---
$ws.Range("E10").Value = $false
$ws.Range("F10").Value = "DOES NOT LOAD"
$ws.Range("G10").Value = 2021
$ws.Range("H10").Value = "tried using multiple different browsers and devices, website would not load"

# --- Row 11 (MBDC1074 / F-SNP) ---
$ws.Range("E11").Value = $false
$ws.Range("F11").Value = "DOES NOT LOAD"
$ws.Range("G11").Value = 2018
$ws.Range("H11").Value = "tried using multiple different browsers and devices, website would not load"

# --- Row 19 (MBDC1486 / RNA CoSSMos) ---
$ws.Range("G19").Value = "april 18 2022"
$ws.Range("F19").Value = "Cannot find server"
$ws.Range("H19").Value = "wayback machine says it has not archived the URL for the most recent snapshot, does not show site"

# --- Row 26 (MBDC1686 / Hemolytik) ---
$ws.Range("D26").Value = $false

# --- Selection moves to H20 ---
$ws.Range("H20").Select() | Out-Null
